$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 8).Value = 3.8  # H2: 3.9 -> 3.8
$ws.Cells.Item(2, 10).Value = 2.3  # J2: 2.25 -> 2.3
$ws.Cells.Item(2, 11).Value = 2.25  # K2: 2.3 -> 2.25
$ws.Cells.Item(2, 12).Value = 5  # L2: 4.75 -> 5
$ws.Cells.Item(2, 21).Value = 1.36  # U2: 1.33 -> 1.36
$ws.Cells.Item(2, 22).Value = 3  # V2: 3.25 -> 3
$ws.Cells.Item(2, 23).Value = 1.75  # W2: 1.7 -> 1.75
$ws.Cells.Item(2, 24).Value = 2  # X2: 2.05 -> 2
$ws.Cells.Item(2, 25).Value = 8.5  # Y2: 9 -> 8.5
$ws.Cells.Item(2, 26).Value = 9  # Z2: 9.5 -> 9
$ws.Cells.Item(2, 31).Value = 12  # AE2: 13 -> 12
$ws.Cells.Item(2, 36).Value = 15  # AJ2: 17 -> 15
$ws.Cells.Item(3, 15).Value = 1.36  # O3: 1.33 -> 1.36
$ws.Cells.Item(3, 16).Value = 3.2  # P3: 3.4 -> 3.2
$ws.Cells.Item(4, 7).Value = 1.27  # G4: 1.3 -> 1.27
$ws.Cells.Item(4, 8).Value = 5  # H4: 4.5 -> 5
$ws.Cells.Item(4, 9).Value = 13  # I4: 12 -> 13
$ws.Cells.Item(4, 10).Value = 1.73  # J4: 1.8 -> 1.73
$ws.Cells.Item(4, 11).Value = 2.4  # K4: 2.38 -> 2.4
$ws.Cells.Item(4, 12).Value = 10  # L4: 9.5 -> 10
$ws.Cells.Item(4, 15).Value = 1.25  # O4: 1.29 -> 1.25
$ws.Cells.Item(4, 16).Value = 3.75  # P4: 3.5 -> 3.75
$ws.Cells.Item(4, 17).Value = 1.9  # Q4: 1.93 -> 1.9
$ws.Cells.Item(4, 18).Value = 1.95  # R4: 1.93 -> 1.95
$ws.Cells.Item(4, 21).Value = 1.36  # U4: 1.4 -> 1.36
$ws.Cells.Item(4, 22).Value = 3  # V4: 2.75 -> 3
$ws.Cells.Item(4, 25).Value = 6  # Y4: 5.5 -> 6
$ws.Cells.Item(4, 28).Value = 7  # AB4: 7.5 -> 7
$ws.Cells.Item(4, 31).Value = 10  # AE4: 9 -> 10
$ws.Cells.Item(4, 32).Value = 10  # AF4: 9 -> 10
$ws.Cells.Item(4, 33).Value = 29  # AG4: 26 -> 29
$ws.Cells.Item(4, 41).Value = 101  # AO4: 81 -> 101
$ws.Cells.Item(4, 44).Value = 2.48  # AR4: 2.65 -> 2.48
$ws.Cells.Item(4, 45).Value = 1.52  # AS4: 1.47 -> 1.52
$ws.Cells.Item(5, 7).Value = 2.35  # G5: 2.3 -> 2.35
$ws.Cells.Item(5, 8).Value = 3.4  # H5: 3.5 -> 3.4
$ws.Cells.Item(5, 9).Value = 2.9  # I5: 2.8 -> 2.9
$ws.Cells.Item(5, 10).Value = 3  # J5: 2.88 -> 3
$ws.Cells.Item(5, 11).Value = 2.25  # K5: 2.3 -> 2.25
$ws.Cells.Item(5, 12).Value = 3.4  # L5: 3.25 -> 3.4
$ws.Cells.Item(5, 13).Value = 1.04  # M5: 1.03 -> 1.04
$ws.Cells.Item(5, 14).Value = 13  # N5: 15 -> 13
$ws.Cells.Item(5, 15).Value = 1.22  # O5: 1.17 -> 1.22
$ws.Cells.Item(5, 16).Value = 4  # P5: 5 -> 4
$ws.Cells.Item(5, 17).Value = 1.75  # Q5: 1.6 -> 1.75
$ws.Cells.Item(5, 18).Value = 2.05  # R5: 2.3 -> 2.05
$ws.Cells.Item(5, 19).Value = 2.75  # S5: 2.38 -> 2.75
$ws.Cells.Item(5, 20).Value = 1.4  # T5: 1.53 -> 1.4
$ws.Cells.Item(5, 21).Value = 1.33  # U5: 1.3 -> 1.33
$ws.Cells.Item(5, 22).Value = 3.25  # V5: 3.4 -> 3.25
$ws.Cells.Item(5, 23).Value = 1.57  # W5: 1.5 -> 1.57
$ws.Cells.Item(5, 24).Value = 2.25  # X5: 2.5 -> 2.25
$ws.Cells.Item(5, 25).Value = 10  # Y5: 11 -> 10
$ws.Cells.Item(5, 30).Value = 23  # AD5: 21 -> 23
$ws.Cells.Item(5, 31).Value = 13  # AE5: 15 -> 13
$ws.Cells.Item(5, 32).Value = 6.5  # AF5: 7 -> 6.5
$ws.Cells.Item(5, 33).Value = 12  # AG5: 11 -> 12
$ws.Cells.Item(5, 34).Value = 41  # AH5: 34 -> 41
$ws.Cells.Item(5, 35).Value = 151  # AI5: 101 -> 151
$ws.Cells.Item(5, 36).Value = 11  # AJ5: 13 -> 11
$ws.Cells.Item(5, 37).Value = 15  # AK5: 17 -> 15
$ws.Cells.Item(5, 41).Value = 26  # AO5: 23 -> 26
$ws.Cells.Item(5, 44).Value = 2.03  # AR5: 1.98 -> 2.03
$ws.Cells.Item(5, 45).Value = 1.83  # AS5: 1.88 -> 1.83
$ws.Cells.Item(8, 7).Value = 2.22  # G8: 2.15 -> 2.22
$ws.Cells.Item(8, 9).Value = 3.55  # I8: 3.7 -> 3.55
$ws.Cells.Item(8, 10).Value = 2.9  # J8: 2.82 -> 2.9
$ws.Cells.Item(8, 12).Value = 4.1  # L8: 4.15 -> 4.1
$ws.Cells.Item(8, 14).Value = 5.5  # N8: 5.6 -> 5.5
$ws.Cells.Item(8, 16).Value = 2.5  # P8: 2.52 -> 2.5
$ws.Cells.Item(8, 17).Value = 2.37  # Q8: 2.35 -> 2.37
$ws.Cells.Item(8, 18).Value = 1.52  # R8: 1.53 -> 1.52
$ws.Cells.Item(8, 19).Value = 4.25  # S8: 4.15 -> 4.25
$ws.Cells.Item(8, 20).Value = 1.18  # T8: 1.19 -> 1.18
$ws.Cells.Item(8, 21).Value = 1.52  # U8: 1.5 -> 1.52
$ws.Cells.Item(8, 22).Value = 2.37  # V8: 2.42 -> 2.37
$ws.Cells.Item(8, 23).Value = 2  # W8: 1.98 -> 2
$ws.Cells.Item(8, 24).Value = 1.72  # X8: 1.75 -> 1.72
$ws.Cells.Item(8, 25).Value = 5.9  # Y8: 5.8 -> 5.9
$ws.Cells.Item(8, 26).Value = 9.5  # Z8: 9.25 -> 9.5
$ws.Cells.Item(8, 27).Value = 9.25  # AA8: 9 -> 9.25
$ws.Cells.Item(8, 28).Value = 22  # AB8: 21 -> 22
$ws.Cells.Item(8, 30).Value = 40  # AD8: 37 -> 40
$ws.Cells.Item(8, 31).Value = 5.5  # AE8: 5.6 -> 5.5
$ws.Cells.Item(8, 33).Value = 16  # AG8: 15.5 -> 16
$ws.Cells.Item(8, 34).Value = 100  # AH8: 90 -> 100
$ws.Cells.Item(8, 36).Value = 8.25  # AJ8: 9 -> 8.25
$ws.Cells.Item(8, 37).Value = 17.5  # AK8: 19.5 -> 17.5
$ws.Cells.Item(8, 39).Value = 55  # AM8: 60 -> 55
$ws.Cells.Item(8, 41).Value = 50  # AO8: 45 -> 50
$ws.Cells.Item(15, 7).Value = 4.5  # G15: 4.2 -> 4.5
$ws.Cells.Item(15, 8).Value = 3.25  # H15: 3.2 -> 3.25
$ws.Cells.Item(15, 9).Value = 1.75  # I15: 1.85 -> 1.75
$ws.Cells.Item(15, 10).Value = 5.5  # J15: 5 -> 5.5
$ws.Cells.Item(15, 12).Value = 2.5  # L15: 2.63 -> 2.5
$ws.Cells.Item(15, 23).Value = 2.2  # W15: 2.1 -> 2.2
$ws.Cells.Item(15, 24).Value = 1.62  # X15: 1.67 -> 1.62
$ws.Cells.Item(15, 25).Value = 10  # Y15: 9.5 -> 10
$ws.Cells.Item(15, 26).Value = 23  # Z15: 21 -> 23
$ws.Cells.Item(15, 27).Value = 17  # AA15: 15 -> 17
$ws.Cells.Item(15, 31).Value = 7  # AE15: 7.5 -> 7
$ws.Cells.Item(15, 33).Value = 21  # AG15: 19 -> 21
$ws.Cells.Item(15, 34).Value = 81  # AH15: 67 -> 81
$ws.Cells.Item(15, 36).Value = 5.5  # AJ15: 6 -> 5.5
$ws.Cells.Item(15, 37).Value = 7.5  # AK15: 8 -> 7.5
$ws.Cells.Item(15, 39).Value = 13  # AM15: 15 -> 13
$ws.Cells.Item(17, 8).Value = 2.9  # H17: 3.1 -> 2.9
$ws.Cells.Item(17, 9).Value = 3.25  # I17: 3.1 -> 3.25
$ws.Cells.Item(17, 11).Value = 1.95  # K17: 2 -> 1.95
$ws.Cells.Item(17, 13).Value = 1.1  # M17: 1.08 -> 1.1
$ws.Cells.Item(17, 14).Value = 7  # N17: 7.5 -> 7
$ws.Cells.Item(17, 19).Value = 4.5  # S17: 4.33 -> 4.5
$ws.Cells.Item(17, 20).Value = 1.18  # T17: 1.2 -> 1.18
$ws.Cells.Item(17, 31).Value = 7  # AE17: 7.5 -> 7
$ws.Cells.Item(17, 35).Value = 800  # AI17: 900 -> 800
$ws.Cells.Item(17, 36).Value = 8.5  # AJ17: 8 -> 8.5
$ws.Cells.Item(17, 38).Value = 13  # AL17: 12 -> 13
$ws.Cells.Item(17, 42).Value = 1.78  # AP17: 1.75 -> 1.78
$ws.Cells.Item(17, 43).Value = 2.03  # AQ17: 2.05 -> 2.03
$ws.Cells.Item(18, 17).Value = 1.88  # Q18: 1.9 -> 1.88
$ws.Cells.Item(18, 18).Value = 1.98  # R18: 1.95 -> 1.98
$ws.Cells.Item(23, 7).Value = 1.62  # G23: 1.57 -> 1.62
$ws.Cells.Item(23, 8).Value = 3.85  # H23: 3.75 -> 3.85
$ws.Cells.Item(23, 9).Value = 4.5  # I23: 5.1 -> 4.5
$ws.Cells.Item(23, 10).Value = 2.18  # J23: 2.15 -> 2.18
$ws.Cells.Item(23, 11).Value = 2.27  # K23: 2.22 -> 2.27
$ws.Cells.Item(23, 12).Value = 4.7  # L23: 5.2 -> 4.7
$ws.Cells.Item(23, 19).Value = 2.62  # S23: 2.6 -> 2.62
$ws.Cells.Item(23, 21).Value = 1.35  # U23: 1.37 -> 1.35
$ws.Cells.Item(23, 22).Value = 2.95  # V23: 2.85 -> 2.95
$ws.Cells.Item(23, 23).Value = 1.72  # W23: 1.7 -> 1.72
$ws.Cells.Item(23, 27).Value = 8  # AA23: 7.9 -> 8
$ws.Cells.Item(23, 28).Value = 12.5  # AB23: 12 -> 12.5
$ws.Cells.Item(23, 29).Value = 12.5  # AC23: 11.75 -> 12.5
$ws.Cells.Item(23, 30).Value = 23  # AD23: 22 -> 23
$ws.Cells.Item(23, 32).Value = 7.8  # AF23: 7.6 -> 7.8
$ws.Cells.Item(23, 36).Value = 15  # AJ23: 15.5 -> 15
$ws.Cells.Item(23, 37).Value = 28  # AK23: 32 -> 28
$ws.Cells.Item(23, 38).Value = 15  # AL23: 16.5 -> 15
$ws.Cells.Item(23, 39).Value = 75  # AM23: 100 -> 75
$ws.Cells.Item(23, 40).Value = 40  # AN23: 50 -> 40
$ws.Cells.Item(23, 41).Value = 40  # AO23: 45 -> 40
$ws.Cells.Item(24, 7).Value = 1.93  # G24: 2.3 -> 1.93
$ws.Cells.Item(24, 8).Value = 3.55  # H24: 3.25 -> 3.55
$ws.Cells.Item(24, 9).Value = 3.3  # I24: 2.8 -> 3.3
$ws.Cells.Item(24, 10).Value = 2.5  # J24: 2.9 -> 2.5
$ws.Cells.Item(24, 11).Value = 2.25  # K24: 2.12 -> 2.25
$ws.Cells.Item(24, 12).Value = 3.8  # L24: 3.4 -> 3.8
$ws.Cells.Item(24, 14).Value = 8  # N24: 7.6 -> 8
$ws.Cells.Item(24, 15).Value = 1.25  # O24: 1.27 -> 1.25
$ws.Cells.Item(24, 16).Value = 3.6  # P24: 3.45 -> 3.6
$ws.Cells.Item(24, 17).Value = 1.72  # Q24: 1.8 -> 1.72
$ws.Cells.Item(24, 18).Value = 2  # R24: 1.91 -> 2
$ws.Cells.Item(24, 19).Value = 2.75  # S24: 2.9 -> 2.75
$ws.Cells.Item(24, 20).Value = 1.4  # T24: 1.36 -> 1.4
$ws.Cells.Item(24, 21).Value = 1.34  # U24: 1.38 -> 1.34
$ws.Cells.Item(24, 22).Value = 3  # V24: 2.82 -> 3
$ws.Cells.Item(24, 24).Value = 2.1  # X24: 2.12 -> 2.1
$ws.Cells.Item(24, 25).Value = 8.5  # Y24: 8.75 -> 8.5
$ws.Cells.Item(24, 26).Value = 10  # Z24: 12 -> 10
$ws.Cells.Item(24, 27).Value = 8.5  # AA24: 9 -> 8.5
$ws.Cells.Item(24, 28).Value = 17.5  # AB24: 24 -> 17.5
$ws.Cells.Item(24, 29).Value = 14.5  # AC24: 18 -> 14.5
$ws.Cells.Item(24, 30).Value = 23  # AD24: 25 -> 23
$ws.Cells.Item(24, 31).Value = 8  # AE24: 7.6 -> 8
$ws.Cells.Item(24, 32).Value = 7  # AF24: 6.5 -> 7
$ws.Cells.Item(24, 33).Value = 13.5  # AG24: 12.5 -> 13.5
$ws.Cells.Item(24, 34).Value = 55  # AH24: 50 -> 55
$ws.Cells.Item(24, 35).Value = 400  # AI24: 350 -> 400
$ws.Cells.Item(24, 36).Value = 11.5  # AJ24: 9.75 -> 11.5
$ws.Cells.Item(24, 37).Value = 18.5  # AK24: 15.5 -> 18.5
$ws.Cells.Item(24, 38).Value = 11.75  # AL24: 10.25 -> 11.75
$ws.Cells.Item(24, 39).Value = 45  # AM24: 35 -> 45
$ws.Cells.Item(24, 40).Value = 28  # AN24: 23 -> 28
$ws.Cells.Item(24, 41).Value = 32  # AO24: 29 -> 32
